$wb = $excel.ActiveWorkbook

# Rename existing sheet from Sheet1 to cellVoltage
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "cellVoltage"

# Add a new worksheet for ADC current testing, placed right after cellVoltage
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "currentADC"

# Populate header row + first data value
$ws2.Range("A1").Value = "Current"
$ws2.Range("B1").Value = "ADC Value"
$ws2.Range("A2").Value = 0

# Make the new sheet the active / selected sheet
$ws2.Activate()
$ws2.Range("A2").Select()
